$wb = $excel.ActiveWorkbook

# --- Sheet "Solucion": shuffle the Pedido/Salida assignment (new randomized run) ---
$wsSol = $wb.Worksheets.Item("Solucion")

$solucionData = @(
    @{Row=2; A="Pedido_5"; B="S001"},
    @{Row=3; A="Pedido_55"; B="S021"},
    @{Row=4; A="Pedido_59"; B="S041"},
    @{Row=5; A="Pedido_67"; B="S061"},
    @{Row=6; A="Pedido_13"; B="S071"},
    @{Row=7; A="Pedido_43"; B="S031"},
    @{Row=8; A="Pedido_57"; B="S051"},
    @{Row=9; A="Pedido_12"; B="S011"},
    @{Row=10; A="Pedido_30"; B="S042"},
    @{Row=11; A="Pedido_29"; B="S022"},
    @{Row=12; A="Pedido_64"; B="S062"},
    @{Row=13; A="Pedido_71"; B="S052"},
    @{Row=14; A="Pedido_33"; B="S002"},
    @{Row=15; A="Pedido_25"; B="S032"},
    @{Row=16; A="Pedido_77"; B="S072"},
    @{Row=17; A="Pedido_63"; B="S043"},
    @{Row=18; A="Pedido_3"; B="S063"},
    @{Row=19; A="Pedido_75"; B="S023"},
    @{Row=20; A="Pedido_28"; B="S053"},
    @{Row=21; A="Pedido_65"; B="S012"},
    @{Row=22; A="Pedido_22"; B="S033"},
    @{Row=23; A="Pedido_37"; B="S003"},
    @{Row=24; A="Pedido_73"; B="S044"},
    @{Row=25; A="Pedido_50"; B="S024"},
    @{Row=26; A="Pedido_47"; B="S073"},
    @{Row=27; A="Pedido_34"; B="S013"},
    @{Row=28; A="Pedido_35"; B="S054"},
    @{Row=29; A="Pedido_76"; B="S034"},
    @{Row=30; A="Pedido_58"; B="S064"},
    @{Row=31; A="Pedido_60"; B="S025"},
    @{Row=32; A="Pedido_49"; B="S074"},
    @{Row=33; A="Pedido_42"; B="S045"},
    @{Row=34; A="Pedido_24"; B="S004"},
    @{Row=35; A="Pedido_74"; B="S055"},
    @{Row=36; A="Pedido_14"; B="S035"},
    @{Row=37; A="Pedido_79"; B="S065"},
    @{Row=38; A="Pedido_32"; B="S014"},
    @{Row=39; A="Pedido_46"; B="S026"},
    @{Row=40; A="Pedido_38"; B="S075"},
    @{Row=41; A="Pedido_44"; B="S046"},
    @{Row=42; A="Pedido_27"; B="S005"},
    @{Row=43; A="Pedido_53"; B="S066"},
    @{Row=44; A="Pedido_39"; B="S036"},
    @{Row=45; A="Pedido_45"; B="S056"},
    @{Row=46; A="Pedido_68"; B="S015"},
    @{Row=47; A="Pedido_1"; B="S006"},
    @{Row=48; A="Pedido_40"; B="S076"},
    @{Row=49; A="Pedido_80"; B="S027"},
    @{Row=50; A="Pedido_54"; B="S047"},
    @{Row=51; A="Pedido_9"; B="S016"},
    @{Row=52; A="Pedido_56"; B="S067"},
    @{Row=53; A="Pedido_19"; B="S037"},
    @{Row=54; A="Pedido_23"; B="S057"},
    @{Row=55; A="Pedido_31"; B="S007"},
    @{Row=56; A="Pedido_11"; B="S077"},
    @{Row=57; A="Pedido_20"; B="S028"},
    @{Row=58; A="Pedido_21"; B="S048"},
    @{Row=59; A="Pedido_7"; B="S017"},
    @{Row=60; A="Pedido_4"; B="S068"},
    @{Row=61; A="Pedido_51"; B="S038"},
    @{Row=62; A="Pedido_36"; B="S058"},
    @{Row=63; A="Pedido_8"; B="S008"},
    @{Row=64; A="Pedido_17"; B="S078"},
    @{Row=65; A="Pedido_26"; B="S029"},
    @{Row=66; A="Pedido_72"; B="S049"},
    @{Row=67; A="Pedido_62"; B="S018"},
    @{Row=68; A="Pedido_15"; B="S059"},
    @{Row=69; A="Pedido_70"; B="S069"},
    @{Row=70; A="Pedido_66"; B="S009"},
    @{Row=71; A="Pedido_2"; B="S039"},
    @{Row=72; A="Pedido_52"; B="S079"},
    @{Row=73; A="Pedido_41"; B="S050"},
    @{Row=74; A="Pedido_48"; B="S019"},
    @{Row=75; A="Pedido_69"; B="S010"},
    @{Row=76; A="Pedido_18"; B="S030"},
    @{Row=77; A="Pedido_6"; B="S070"},
    @{Row=78; A="Pedido_10"; B="S060"},
    @{Row=79; A="Pedido_78"; B="S020"},
    @{Row=80; A="Pedido_61"; B="S040"},
    @{Row=81; A="Pedido_16"; B="S080"}
)

foreach ($row in $solucionData) {
    $wsSol.Cells.Item($row.Row, 1).Value = $row.A
    $wsSol.Cells.Item($row.Row, 2).Value = $row.B
}

# --- Sheet "Metricas": updated per-zone Tiempo values from the new run ---
$wsMet = $wb.Worksheets.Item("Metricas")
$wsMet.Range("B2").Value = 536.8818845280571
$wsMet.Range("B3").Value = 540.553243593902
$wsMet.Range("B4").Value = 530.5490539517785
$wsMet.Range("B5").Value = 532.9112336468809

# --- Sheet "Resumen": the zone with the maximum time, and that maximum ---
$wsRes = $wb.Worksheets.Item("Resumen")
$wsRes.Range("B2").Value = "Z2"
$wsRes.Range("C2").Value = 540.553243593902
